# Natmi following Dr Hou advice
# Rebuild the LR-pair data block (rows 2-10) on Sheet1 for Pdgfb-Lrp1 signalling,
# expanding the 2x... triad of sending/target clusters (ECs/FAPs/sCs) into a full 3x3
# cross-product block with refreshed NATMI statistics.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$natmiRows = @(
  @("ECs", "Pdgfb", "Lrp1", "ECs", 3, 1, 41.23709466666667, 123.711284, 0.956365997213294, 0.9563659972132939, 3, 1, 17.16653, 51.49959, 0.0560345397128279, 0.0560345397128279, 707.8978227081733, 6371.08040437356, 0.05358952845084657, 0.05358952845084657),
  @("ECs", "Pdgfb", "Lrp1", "FAPs", 3, 1, 41.23709466666667, 123.711284, 0.956365997213294, 0.9563659972132939, 3, 1, 256.4443053333333, 769.332916, 0.8370788162388805, 0.8370788162388805, 10575.01809575824, 95175.16286182414, 0.8005537168384206, 0.8005537168384205),
  @("ECs", "Pdgfb", "Lrp1", "sCs", 3, 1, 41.23709466666667, 123.711284, 0.956365997213294, 0.9563659972132939, 3, 1, 32.74538866666666, 98.236166, 0.1068866440482915, 0.1068866440482915, 1350.324692344127, 12152.92223109714, 0.1022227519240267, 0.1022227519240267),
  @("FAPs", "Pdgfb", "Lrp1", "ECs", 1, 0.3333333333333333, 0.030543, 0.091629, 0.0007083497731593903, 0.0007083497731593902, 3, 1, 17.16653, 51.49959, 0.0560345397128279, 0.0560345397128279, 0.5243173257899999, 4.71885593211, 0.00003969205349467248, 0.00003969205349467248),
  @("FAPs", "Pdgfb", "Lrp1", "FAPs", 1, 0.3333333333333333, 0.030543, 0.091629, 0.0007083497731593903, 0.0007083497731593902, 3, 1, 256.4443053333333, 769.332916, 0.8370788162388805, 0.8370788162388805, 7.832578417795999, 70.493205760164, 0.000592944589599342, 0.0005929445895993419),
  @("FAPs", "Pdgfb", "Lrp1", "sCs", 1, 0.3333333333333333, 0.030543, 0.091629, 0.0007083497731593903, 0.0007083497731593902, 3, 1, 32.74538866666666, 98.236166, 0.1068866440482915, 0.1068866440482915, 1.000142406046, 9.001281654413999, 0.0000757131300653758, 0.00007571313006537579),
  @("sCs", "Pdgfb", "Lrp1", "ECs", 3, 1, 1.850891, 5.552673, 0.0429256530135467, 0.04292565301354669, 3, 1, 17.16653, 51.49959, 0.0560345397128279, 0.0560345397128279, 31.77337587823, 285.9603829040699, 0.002405319208486653, 0.002405319208486652),
  @("sCs", "Pdgfb", "Lrp1", "FAPs", 3, 1, 1.850891, 5.552673, 0.0429256530135467, 0.04292565301354669, 3, 1, 256.4443053333333, 769.332916, 0.8370788162388805, 0.8370788162388805, 474.6504567427186, 4271.854110684468, 0.0359321548108606, 0.03593215481086059),
  @("sCs", "Pdgfb", "Lrp1", "sCs", 3, 1, 1.850891, 5.552673, 0.0429256530135467, 0.04292565301354669, 3, 1, 32.74538866666666, 98.236166, 0.1068866440482915, 0.1068866440482915, 60.60814517463533, 545.4733065717179, 0.004588178994199439, 0.004588178994199438),
)

$r = 2
foreach ($rowValues in $natmiRows) {
  $c = 1
  foreach ($val in $rowValues) {
    $ws.Cells.Item($r, $c).Value = $val
    $c++
  }
  $r++
}

